# spec sheet for macro specifications
# Adds a "macro type" / package (LEF/LIB) block to the top of the
# ADC_100MS spec sheet, just above the existing "width"/"height" rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADC_100MS")

# Make room for two new rows above the "width"/"height" pair (old rows
# 3-4 stay put; everything from the old row 8 header downward shifts
# down by two, carrying its formulas/styles with it).
$ws.Rows("5:6").Insert()

# New "macro type" row, right above "width".
$ws.Range("B2").Value = "macro type"
$ws.Range("C2").Value = "mixed signal"

# New LEF/LIB rows, right below "height".
$ws.Range("B5").Value = "LEF"
$ws.Range("C5").Value = "test.lef"
$ws.Range("B6").Value = "LIB"
$ws.Range("C6").Value = "test.lib"

# The value column of this little key/value block is right aligned.
$ws.Range("C3:C6").HorizontalAlignment = -4152

# Park the selection where the author left it.
$ws.Range("D3").Select()
